# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# figures for each coin row. Price cells that look like plain numbers
# ("1.007", "314.52", ...) are forced to text first (NumberFormat "@")
# so Excel keeps them as literal strings instead of parsing them as
# numeric values, then the format is reset back to Normal so no stray
# cell style is left behind. Volume percentages already contain spaces
# and a "%" sign so they stay text without any extra handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.216.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.823.05'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.99%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.18%  '
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4268'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.70%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3684'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07240'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.00'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.826.52'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.725'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07089'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.314'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('E20').Value = '  -3.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.249.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.139'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.052.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.005'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.130'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.231'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.32'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08894'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.196'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7560'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.436'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.839'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.88%  '
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.113'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05271'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.128'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.871'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1696'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5038'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '107.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4735'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06371'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('E50').Value = '  -3.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.808'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.27%  '
